$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Change 1: insert a new "17 May 1778" entry before "20 May 1778"
# ---------------------------------------------------------------
$find1 = $d.Content
$found1 = $find1.Find.Execute("20 May 1778", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found1) {
    $insertPoint = $d.Range($find1.Start, $find1.Start)
    $insertPoint.InsertParagraphBefore()

    # locate the freshly-created (now empty) paragraph that sits right
    # before the "20 May 1778" paragraph
    $newParaStart = $find1.Start
    # After InsertParagraphBefore, the empty new paragraph begins where
    # "20 May 1778" used to begin.
    $dateText = "17 May 1778"
    $dateRange = $d.Range($newParaStart, $newParaStart)
    $dateRange.InsertAfter($dateText)

    $dateEnd = $newParaStart + $dateText.Length
    $boldRange = $d.Range($newParaStart, $dateEnd)
    $boldRange.Font.Bold = 1
    $boldRange.Font.Color = 0

    $descText = "  Christian Gottlob Neefe (30) marries the singer and actress Suzanne Zinck in Frankfurt.  She is the adopted daughter of Georg Benda (55)."
    $descInsertPoint = $d.Range($dateEnd, $dateEnd)
    $descInsertPoint.InsertAfter($descText)

    $descRange = $d.Range($dateEnd, $dateEnd + $descText.Length)
    $descRange.Font.Bold = 0
    $descRange.Font.Color = 0
}

# ---------------------------------------------------------------
# Change 2: "January 2016" -> "May" + " 2016" (two separate runs)
# ---------------------------------------------------------------
$find2 = $d.Content
$found2 = $find2.Find.Execute("January 2016", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found2) {
    $startPos = $find2.Start
    $find2.Delete()

    $mayInsert = $d.Range($startPos, $startPos)
    $mayInsert.InsertAfter("May")

    $afterMay = $startPos + 3
    $yearInsert = $d.Range($afterMay, $afterMay)
    $yearInsert.InsertAfter(" 2016")
}
